# 3.9.2 "Mortality rate attributed to unsafe water..." indicator sheet
# update:
#   1) Fix the Kyrgyz-language title in A1 (typo correction:
#      "Коопсуз...жоктугунана" -> "Коопсуздук...жоктугунан").
#   2) Add a new "2022" data column (after the existing 2021 column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Corrected Kyrgyz title in A1.
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# 2) Add the 2022 column. Duplicating column R (2021) into the new column S
#    first gives the new column the same number formatting/border/font
#    styling as the rest of the year columns; the actual 2022 figures are
#    then written on top of the copied 2021 figures below.
$ws.Columns("R").Copy()
$ws.Columns("S").Insert()

$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 1.2
$ws.Range("S6").Value = 2.7
$ws.Range("S7").Value = 0.9
$ws.Range("S8").Value = 0.4
$ws.Range("S9").Value = 0.7
$ws.Range("S10").Value = 0.9
$ws.Range("S11").Value = 1.1
$ws.Range("S12").Value = 2.7
$ws.Range("S13").Value = 0.4
$ws.Range("S14").Value = 0.6
